# Weight & Balance recalculation update (ACBalanceManager / CabinConfiguration work):
# the single-passenger design mass dropped 106 -> 105 kg, which ripples through
# every downstream weight-estimation-method table across all sheets, and several
# per-component "methods comparison" tables were re-ordered/re-computed as well.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("GLOBAL RESULTS")
$ws.Range("C3").Value = 105.0
$ws.Range("C5").Value = 67438.34854338123
$ws.Range("C6").Value = 67438.34854338123
$ws.Range("C7").Value = 65415.198087079785
$ws.Range("C9").Value = 18772.130374844557
$ws.Range("C10").Value = 48666.21816853667
$ws.Range("C11").Value = 48666.21816853667
$ws.Range("C12").Value = 13650.0
$ws.Range("C13").Value = 13650.0
$ws.Range("C14").Value = 35016.218168536674
$ws.Range("C15").Value = 34296.435556153665
$ws.Range("C16").Value = 337.2098698830023
$ws.Range("C19").Value = 32687.545426036682
$ws.Range("C20").Value = 18283.82117645248

$ws = $wb.Worksheets("FUSELAGE")
# method rows (8-15) keep JENKINSON/ROSKAM/KROO/TORENBEEK_2013/NICOLAI_1984/RAYMER/TORENBEEK_1976/SADRAEY order; values recomputed
$ws.Range("C2").Value = 6744.197397660046
$ws.Range("C3").Value = 7247.125
$ws.Range("D3").Value = 7.45718982831743
$ws.Range("C5").Value = 7247.124999999999
$ws.Range("C8").Value = 6645.0
$ws.Range("D8").Value = -1.4708554897053137
$ws.Range("C9").Value = 6931.0
$ws.Range("D9").Value = 2.7698270279687693
$ws.Range("C10").Value = 6493.0
$ws.Range("D10").Value = -3.724644799797833
$ws.Range("D11").Value = 17.004285828553105
$ws.Range("C12").Value = 9233.0
$ws.Range("D12").Value = 36.902873026869955
$ws.Range("C13").Value = 8148.0
$ws.Range("D13").Value = 20.814969069959535
$ws.Range("C14").Value = 6240.0
$ws.Range("D14").Value = -7.4760177962018295
$ws.Range("D15").Value = -5.162918241106875

$ws = $wb.Worksheets("WING")
# method rows 8-12 re-ordered: TORENBEEK_1982->RAYMER, ROSKAM->KROO, JENKINSON->TORENBEEK_2013, KROO->TORENBEEK_1982, TORENBEEK_2013->JENKINSON
$ws.Range("C2").Value = 7148.849241519648
$ws.Range("C3").Value = 5479.714285714284
$ws.Range("D3").Value = -23.34830263465663
$ws.Range("C5").Value = 5479.714285714284
$ws.Range("A8").Value = "RAYMER"
$ws.Range("C8").Value = 8372.0
$ws.Range("D8").Value = 17.10975734914708
$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 7124.0
$ws.Range("D9").Value = -0.3475977836450319
$ws.Range("A10").Value = "TORENBEEK_2013"
$ws.Range("C10").Value = 5858.0
$ws.Range("D10").Value = -18.056741692390876
$ws.Range("A11").Value = "TORENBEEK_1982"
$ws.Range("C11").Value = 6037.0
$ws.Range("D11").Value = -15.552842198184313
$ws.Range("A12").Value = "JENKINSON"
$ws.Range("C12").Value = 886.0
$ws.Range("D12").Value = -87.60639691694406
$ws.Range("C13").Value = 4280.0
$ws.Range("D13").Value = -40.13022438433475
$ws.Range("D14").Value = -18.85407281624436

$ws = $wb.Worksheets("HORIZONTAL TAIL")
# method rows 8-12 re-ordered: NICOLAI_2013->RAYMER, JENKINSON->NICOLAI_2013, ROSKAM->KROO, KROO->HOWE, HOWE->JENKINSON
$ws.Range("C2").Value = 775.5827007309051
$ws.Range("C3").Value = 624.5714285714284
$ws.Range("D3").Value = -19.470685978060676
$ws.Range("C5").Value = 624.5714285714284
$ws.Range("A8").Value = "RAYMER"
$ws.Range("C8").Value = 525.0
$ws.Range("D8").Value = -32.308959508090794
$ws.Range("A9").Value = "NICOLAI_2013"
$ws.Range("C9").Value = 415.0
$ws.Range("D9").Value = -46.49184418258606
$ws.Range("A10").Value = "KROO"
$ws.Range("C10").Value = 738.0
$ws.Range("D10").Value = -4.845737365659056
$ws.Range("A11").Value = "HOWE"
$ws.Range("C11").Value = 472.0
$ws.Range("D11").Value = -39.142531214893054
$ws.Range("A12").Value = "JENKINSON"
$ws.Range("C12").Value = 700.0
$ws.Range("D12").Value = -9.745279344121055
$ws.Range("C13").Value = 482.0
$ws.Range("D13").Value = -37.85317806266622
$ws.Range("C14").Value = 1040.0
$ws.Range("D14").Value = 34.092727831591574

$ws = $wb.Worksheets("VERTICAL TAIL")
# method rows 9-11 re-ordered: ROSKAM->KROO, KROO->HOWE, HOWE->JENKINSON (row 8 JENKINSON stays put)
$ws.Range("C2").Value = 775.5827007309051
$ws.Range("C3").Value = 470.33333333333326
$ws.Range("D3").Value = -39.35742340693086
$ws.Range("C5").Value = 470.33333333333326
$ws.Range("C8").Value = 194.0
$ws.Range("D8").Value = -74.98654884679927
$ws.Range("A9").Value = "KROO"
$ws.Range("C9").Value = 497.0
$ws.Range("D9").Value = -35.91914833432595
$ws.Range("A10").Value = "HOWE"
$ws.Range("C10").Value = 382.0
$ws.Range("D10").Value = -50.746709584934635
$ws.Range("A11").Value = "JENKINSON"
$ws.Range("C11").Value = 502.0
$ws.Range("D11").Value = -35.27447175821253
$ws.Range("C12").Value = 482.0
$ws.Range("D12").Value = -37.85317806266622
$ws.Range("C13").Value = 765.0
$ws.Range("D13").Value = -1.3644838546465818

$ws = $wb.Worksheets("NACELLES")
# NACELLE 1 (rows 10-12) and NACELLE 2 (rows 17-19) methods re-ordered: JENKINSON->KUNDU, ROSKAM->JENKINSON, KUNDU->ROSKAM
$ws.Range("C2").Value = 1281.3975055554085
$ws.Range("D3").Value = 117.10671262733348
$ws.Range("A10").Value = "KUNDU"
$ws.Range("C10").Value = 1389.0
$ws.Range("D10").Value = 116.79455344311013
$ws.Range("A11").Value = "JENKINSON"
$ws.Range("C11").Value = 1410.0
$ws.Range("D11").Value = 120.07222487745521
$ws.Range("A12").Value = "ROSKAM"
$ws.Range("C12").Value = 1374.0
$ws.Range("D12").Value = 114.45335956143508
$ws.Range("A17").Value = "KUNDU"
$ws.Range("C17").Value = 1389.0
$ws.Range("D17").Value = 116.79455344311013
$ws.Range("A18").Value = "JENKINSON"
$ws.Range("C18").Value = 1410.0
$ws.Range("D18").Value = 120.07222487745521
$ws.Range("A19").Value = "ROSKAM"
$ws.Range("C19").Value = 1374.0
$ws.Range("D19").Value = 114.45335956143508

$ws = $wb.Worksheets("POWER PLANT")
# ENGINE 1 (rows 11-13) and ENGINE 2 (rows 18-20) methods re-ordered: TORENBEEK_2013->KUNDU, TORENBEEK_1976->TORENBEEK_2013, KUNDU->TORENBEEK_1976
$ws.Range("C2").Value = 5597.683840057838
$ws.Range("D3").Value = 15.250048371196929
$ws.Range("A11").Value = "KUNDU"
$ws.Range("C11").Value = 3265.0
$ws.Range("D11").Value = 16.655391525873153
$ws.Range("A12").Value = "TORENBEEK_2013"
$ws.Range("C12").Value = 3458.0
$ws.Range("D12").Value = 23.551100733987553
$ws.Range("A13").Value = "TORENBEEK_1976"
$ws.Range("C13").Value = 2954.0
$ws.Range("D13").Value = 5.543652853730258
$ws.Range("A18").Value = "KUNDU"
$ws.Range("C18").Value = 3265.0
$ws.Range("D18").Value = 16.655391525873153
$ws.Range("A19").Value = "TORENBEEK_2013"
$ws.Range("C19").Value = 3458.0
$ws.Range("D19").Value = 23.551100733987553
$ws.Range("A20").Value = "TORENBEEK_1976"
$ws.Range("C20").Value = 2954.0
$ws.Range("D20").Value = 5.543652853730258

$ws = $wb.Worksheets("LANDING GEARS")
$ws.Range("C2").Value = 2765.1209330406186
$ws.Range("C3").Value = 1680.077128833439
$ws.Range("D3").Value = -39.2403743084766
$ws.Range("C5").Value = 1680.0771288334388
$ws.Range("C9").Value = 1680.077128833439
$ws.Range("D9").Value = -39.2403743084766
$ws.Range("C11").Value = 216.3153282268604
$ws.Range("C13").Value = 1463.7618006065786

$ws = $wb.Worksheets("SYSTEMS")
$ws.Range("C2").Value = 9172.10846081766
$ws.Range("C3").Value = 7952.390916250886
$ws.Range("D3").Value = -13.2981151474308
$ws.Range("C4").Value = 7952.390916250884
$ws.Range("C8").Value = 7952.390916250886
$ws.Range("D8").Value = -13.29811514743079
$ws.Range("C21").Value = 988.2918262859084
$ws.Range("C23").Value = 988.2918262859083
$ws.Range("C26").Value = 490.3269234142957
$ws.Range("C28").Value = 490.32692341429566
$ws.Range("C36").Value = 814.6505130614842
$ws.Range("C38").Value = 814.650513061484
$ws.Range("C41").Value = 3087.864478432748
$ws.Range("C43").Value = 3087.8644784327475
